$wb = $excel.ActiveWorkbook

# Personal_JNT (sheet7): fill C2:C6 with "Y" and update selection
$wsJNT = $wb.Worksheets.Item("Personal_JNT")
$wsJNT.Range("C2:C6").Value = "Y"
$wsJNT.Range("C2:C6").Select()

# Personal_EL (sheet8): move F1 header to blank, add C2 "Y", move selection
$wsEL = $wb.Worksheets.Item("Personal_EL")
$wsEL.Range("F1").ClearContents()
$wsEL.Range("C2").Value = "Y"
$wsEL.Range("F1:F1048576").Select()

# Estimated (sheet1) becomes the active/selected tab
$wsEstimated = $wb.Worksheets.Item("Estimated")
$wsEstimated.Activate()
